$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIC")

# --- Remove cells that must no longer exist (old single-row/col header artifacts) ---
$removeRefs = @("A1", "B2", "C2", "D2", "E2", "F2", "G2", "H2", "I2", "J2", "K2", "L2", "M2", "N2")
foreach ($ref in $removeRefs) {
    $ws.Range($ref).Clear()
}

# --- Set text (row/column index labels), reusing existing shared strings ---
$stringCells = @{
    "B1" = "N2"
    "C1" = "CO2"
    "D1" = "C1"
    "E1" = "C2"
    "F1" = "C3"
    "G1" = "C4"
    "H1" = "C5"
    "I1" = "C6"
    "J1" = "PS-1"
    "K1" = "PS-2"
    "L1" = "PS-3"
    "M1" = "PS-4"
    "N1" = "PS-5"
    "A2" = "N2"
    "A3" = "CO2"
    "A4" = "C1"
    "A5" = "C2"
    "A6" = "C3"
    "A7" = "C4"
    "A8" = "C5"
    "A9" = "C6"
    "A10" = "PS-1"
    "A11" = "PS-2"
    "A12" = "PS-3"
    "A13" = "PS-4"
    "A14" = "PS-5"
}
foreach ($ref in $stringCells.Keys) {
    $ws.Range($ref).HorizontalAlignment = -4108
    $ws.Range($ref).Value = $stringCells[$ref]
}

# --- Set numeric matrix values ---
$numberCells = @{
    "B3" = 0.02
    "B4" = 0.06
    "C4" = 0.12
    "B5" = 0.08
    "C5" = 0.15
    "B6" = 0.08
    "C6" = 0.15
    "B7" = 0.08
    "C7" = 0.15
    "B8" = 0.08
    "C8" = 0.15
    "B9" = 0.08
    "C9" = 0.15
    "B10" = 0.08
    "C10" = 0.15
    "D10" = 0.06
    "B11" = 0.08
    "C11" = 0.15
    "D11" = 0.08
    "B12" = 0.08
    "C12" = 0.15
    "D12" = 0.09
    "B13" = 0.08
    "C13" = 0.15
    "D13" = 0.11
    "E13" = 0
    "F13" = 0
    "G13" = 0
    "H13" = 0
    "I13" = 0
    "J13" = 0
    "K13" = 0
    "L13" = 0
    "B14" = 0.08
    "C14" = 0.15
    "D14" = 0.14000000000000001
    "E14" = 0
    "F14" = 0
    "G14" = 0
    "H14" = 0
    "I14" = 0
    "J14" = 0
    "K14" = 0
    "L14" = 0
    "M14" = 0
}
foreach ($ref in $numberCells.Keys) {
    $ws.Range($ref).HorizontalAlignment = -4108
    $ws.Range($ref).Value = $numberCells[$ref]
}

# --- Cells that must exist (centered style) but stay blank ---
$emptyRefs = @("C3", "E3", "F3", "G3", "H3", "I3", "J3", "K3", "L3", "M3", "D4", "F4", "G4", "H4", "I4", "J4", "K4", "L4", "M4", "N4", "E5", "G5", "H5", "I5", "J5", "K5", "L5", "M5", "N5", "F6", "H6", "I6", "J6", "K6", "L6", "M6", "N6", "G7", "I7", "J7", "K7", "L7", "M7", "N7", "H8", "J8", "K8", "L8", "M8", "N8", "I9", "K9", "L9", "M9", "N9", "J10", "L10", "M10", "N10", "K11", "M11", "N11", "L12", "N12", "M13", "N13", "N14")
foreach ($ref in $emptyRefs) {
    $ws.Range($ref).HorizontalAlignment = -4108
    $ws.Range($ref).Value = ""
}

# --- Update selection to match the saved view state ---
$ws.Range("N11").Select()

# --- Sheet1 (second tab) selection only changed, update it too ---
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("C2:O2").Select()

